# Apply the "selecao via mapa e filtros multiplos" data update to the POINTS sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace MSP21 entry with MGO12 entry -----------------------
$ws.Range("A2").Value2 = "MGO12"
$ws.Range("B2").Value2 = "Av Denise Cristina Da Rocha, Em Frente Braunas "
$ws.Range("C2").Value2 = -19.802656
$ws.Range("D2").Value2 = -44.0040275
$ws.Range("E2").Value2 = "https://i.ibb.co/q5xpNnh/AV-DENISE-CRISTINA-DA-ROCHA-EM-FRENTE-BRAUNAS-19-48-09-7-S-44-00-14-5-W.jpg"
$ws.Range("F2").Value2 = $null
$ws.Range("H2").Value2 = "Ribeirão das Neves"
$ws.Range("J2").Value2 = "Minas Gerais"
$ws.Range("K2").Value2 = "Brasil"
$ws.Range("L2").Value2 = "Outdoor "

# --- Row 3: replace MFD457 entry with MRC10 entry -----------------------
$ws.Range("A3").Value2 = "MRC10"
$ws.Range("B3").Value2 = "Av.Brasilia proximo ao clube Lago Azul Duquesa II."
$ws.Range("C3").Value2 = -19.7768379
$ws.Range("D3").Value2 = -43.8970794
$ws.Range("E3").Value2 = "https://i.ibb.co/KmtBT8P/Av-Brasilia-proximo-ao-clube-Lago-Azul-Duquesa-II.jpg"
$ws.Range("F3").Value2 = "Bonó"
$ws.Range("H3").Value2 = "Santa Luzia"
$ws.Range("J3").Value2 = "Minas Gerais"
$ws.Range("K3").Value2 = "Brasil"
$ws.Range("L3").Value2 = "Outdoor "

# --- Remove the now-obsolete rows 4 and 5 (MEM33 / MEM34) --------------
$ws.Rows("4:5").Delete()
